{"js": "// Replace two-digit multiplication equations per the commit diff.\nconst pairs = [\n  [\"11\u00d757=627\", \"82\u00d789=7298\"],\n  [\"84\u00d743=3612\", \"24\u00d773=1752\"],\n  [\"71\u00d715=1065\", \"22\u00d795=2090\"],\n  [\"96\u00d762=5952\", \"36\u00d722=792\"],\n  [\"49\u00d717=833\", \"53\u00d715=795\"],\n  [\"55\u00d736=1980\", \"67\u00d718=1206\"],\n  [\"58\u00d711=638\", \"46\u00d741=1886\"],\n  [\"48\u00d756=2688\", \"17\u00d713=221\"],\n  [\"71\u00d774=5254\", \"60\u00d770=4200\"],\n  [\"54\u00d749=2646\", \"66\u00d776=5016\"],\n  [\"41\u00d716=656\", \"74\u00d734=2516\"],\n  [\"71\u00d734=2414\", \"35\u00d769=2415\"],\n  [\"42\u00d718=756\", \"55\u00d752=2860\"],\n  [\"44\u00d750=2200\", \"63\u00d789=5607\"],\n  [\"54\u00d739=2106\", \"42\u00d755=2310\"],\n  [\"87\u00d792=8004\", \"98\u00d732=3136\"],\n  [\"31\u00d780=2480\", \"80\u00d745=3600\"],\n  [\"51\u00d765=3315\", \"23\u00d760=1380\"],\n  [\"84\u00d723=1932\", \"24\u00d734=816\"],\n  [\"24\u00d729=696\", \"75\u00d742=3150\"],\n  [\"66\u00d767=4422\", \"95\u00d727=2565\"],\n  [\"31\u00d735=1085\", \"46\u00d758=2668\"],\n  [\"42\u00d765=2730\", \"16\u00d746=736\"],\n  [\"26\u00d794=2444\", \"21\u00d756=1176\"],\n  [\"68\u00d782=5576\", \"82\u00d769=5658\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace two-digit multiplication equations per the commit diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"11\u00d757=627\", \"82\u00d789=7298\"),\n    @(\"84\u00d743=3612\", \"24\u00d773=1752\"),\n    @(\"71\u00d715=1065\", \"22\u00d795=2090\"),\n    @(\"96\u00d762=5952\", \"36\u00d722=792\"),\n    @(\"49\u00d717=833\", \"53\u00d715=795\"),\n    @(\"55\u00d736=1980\", \"67\u00d718=1206\"),\n    @(\"58\u00d711=638\", \"46\u00d741=1886\"),\n    @(\"48\u00d756=2688\", \"17\u00d713=221\"),\n    @(\"71\u00d774=5254\", \"60\u00d770=4200\"),\n    @(\"54\u00d749=2646\", \"66\u00d776=5016\"),\n    @(\"41\u00d716=656\", \"74\u00d734=2516\"),\n    @(\"71\u00d734=2414\", \"35\u00d769=2415\"),\n    @(\"42\u00d718=756\", \"55\u00d752=2860\"),\n    @(\"44\u00d750=2200\", \"63\u00d789=5607\"),\n    @(\"54\u00d739=2106\", \"42\u00d755=2310\"),\n    @(\"87\u00d792=8004\", \"98\u00d732=3136\"),\n    @(\"31\u00d780=2480\", \"80\u00d745=3600\"),\n    @(\"51\u00d765=3315\", \"23\u00d760=1380\"),\n    @(\"84\u00d723=1932\", \"24\u00d734=816\"),\n    @(\"24\u00d729=696\", \"75\u00d742=3150\"),\n    @(\"66\u00d767=4422\", \"95\u00d727=2565\"),\n    @(\"31\u00d735=1085\", \"46\u00d758=2668\"),\n    @(\"42\u00d765=2730\", \"16\u00d746=736\"),\n    @(\"26\u00d794=2444\", \"21\u00d756=1176\"),\n    @(\"68\u00d782=5576\", \"82\u00d769=5658\"),\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n}\n\nWrite-Output \"Replacements complete\"\n"}
